$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K135").Value = 59973.88140000001
$ws.Range("M135").Value = -57438.88140000001
